$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$oldValue = "System, dnasr281@gmail.com"
$newValue = "dnasr281@gmail.com, System"

$usedLastRow = $ws.UsedRange.Rows.Count

for ($r = 1; $r -le $usedLastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $v = $cell.Value2
    if ($v -eq $oldValue) {
        $cell.Value = $newValue
    }
}
